$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "reviews_average"
$ws.Range("F1").Value = "latitude"
$ws.Range("G1").Value = "longitude"
$ws.Range("H1").Value = "is_permanently_closed"
$ws.Range("I1").Value = "gmaps_link"
$ws.Range("J1").Value = "latest_review_date"
$ws.Range("K1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
